$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "Katja"
$ws.Range("B12").Value = "katja@gmail.com"
$ws.Range("C12").Value = 3
$ws.Range("C12").HorizontalAlignment = -4108

$ws.Range("A13").Value = "Tanja"
$ws.Range("B13").Value = "tanja@gmail.com"
$ws.Range("C13").Value = 2

$ws.Range("A14").Value = "Julja"
$ws.Range("B14").Value = "julja@gmail.com"
$ws.Range("C14").Value = -1

$ws.Range("F15").Select()
